$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the original "#00-8" expense entry (row 2) - it was deleted by the user,
# shifting the remaining entries (#00-9, #00-10) up.
$ws.Rows.Item(2).Delete()

# Add new expense entries recorded through the app's UI.
# Row 4: #00-11
$ws.Range("A4").Value = "#00-11"
$ws.Range("B4").Value = 43490.001388888886
$ws.Range("C4").Value = "Anuj Pal"
$ws.Range("D4").Value = "tyw"
$ws.Range("E4").Value = "Food Expense"
$ws.Range("F4").Value = "Attached"
$ws.Range("G4").Value = 5464.0
$ws.Range("H4").Value = "Credit Card"

# Row 5: #00-12
$ws.Range("A5").Value = "#00-12"
$ws.Range("B5").Value = 43491.001388888886
$ws.Range("C5").Value = "Anuj Pal"
$ws.Range("D5").Value = "this is test"
$ws.Range("E5").Value = "Travel Expense"
$ws.Range("F5").Value = "Attached"
$ws.Range("G5").Value = 567.0
$ws.Range("H5").Value = "Credit Card"

# Row 6: #00-13
$ws.Range("A6").Value = "#00-13"
$ws.Range("B6").Value = 43471.00208333333
$ws.Range("C6").Value = "Anuj Pal"
$ws.Range("D6").Value = "sdfghj"
$ws.Range("E6").Value = "Travel Expense"
$ws.Range("F6").Value = "Attached"
$ws.Range("G6").Value = 5678.0
$ws.Range("H6").Value = "Hard Cash"

# Row 7: #00-7 (created by a different user)
$ws.Range("A7").Value = "#00-7"
$ws.Range("B7").Value = 43491.001388888886
$ws.Range("C7").Value = "Pradip Kumar"
$ws.Range("D7").Value = "pradip"
$ws.Range("E7").Value = "Food Expense"
$ws.Range("F7").Value = "Attached"
$ws.Range("G7").Value = 987.0
$ws.Range("H7").Value = "Hard Cash"

# Widen the "Created By" column to fit the new, longer name.
$ws.Columns.Item(3).ColumnWidth = 12.166666666666666
